$d = $word.ActiveDocument

# The empty paragraph right after "Different outlines or colours represent
# active or not?" needs a new sentence about switching Snake to a singleton
# class, and the document's "_GoBack" bookmark needs to move onto this new
# text (it currently sits at the very end of the document).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "") {
        $target = $i
        break
    }
}

$p = $d.Paragraphs.Item($target)
$p.Range.Text = "Switch snake to be a singleton class – will only want the creation of one snake for this iteration of the game"

# Re-fetch the paragraph (its Range identity may have shifted after the text
# assignment) and build a range that excludes the trailing paragraph mark so
# the bookmark wraps only the new sentence, matching the target markup.
$p = $d.Paragraphs.Item($target)
$newRange = $d.Range($p.Range.Start, $p.Range.End - 1)

# Adding a bookmark with the same name moves/redefines it rather than
# duplicating it, so this both creates it here and removes it from its old
# location at the tail of the document.
$d.Bookmarks.Add("_GoBack", $newRange)
